# Update gh-pages output data (苏州-漫展信息)
# Sheet "展览" (sheet index 1) and sheet "全部类型" (sheet index 4) both
# receive the same underlying data refresh: row 2's lowest-price cell (G2)
# becomes "不可售" (not for sale), and several "想去人数" (F column) counts
# are refreshed to newer scraped values.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F3").Value = 14781
$ws1.Range("F4").Value = 18208
$ws1.Range("F6").Value = 99
$ws1.Range("F15").Value = 76
$ws1.Range("F18").Value = 1376
$ws1.Range("F21").Value = 77
$ws1.Range("F22").Value = 222
$ws1.Range("F23").Value = 7534
$ws1.Range("F27").Value = 1198
$ws1.Range("F29").Value = 5911
$ws1.Range("F30").Value = 86
$ws1.Range("F31").Value = 51
$ws1.Range("F35").Value = 5223

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F3").Value = 14781
$ws4.Range("F4").Value = 18208
$ws4.Range("F6").Value = 99
$ws4.Range("F15").Value = 76
$ws4.Range("F18").Value = 1376
$ws4.Range("F22").Value = 77
$ws4.Range("F23").Value = 222
$ws4.Range("F24").Value = 7534
$ws4.Range("F28").Value = 1198
$ws4.Range("F31").Value = 5911
$ws4.Range("F32").Value = 86
$ws4.Range("F33").Value = 51
$ws4.Range("F37").Value = 5223
